# refactor: delete IOException on main function
#
# - Sheet "회원 정보": remove the last row (row 8, all "t"/IOException rows)
#   and change row 7's first cell from "r" to "e".
# - Sheet "자기소개서": replace the placeholder self-intro text with the
#   new text, and narrow column A's width.

$wb = $excel.ActiveWorkbook

$wsMember = $wb.Worksheets.Item(1)
$wsIntro  = $wb.Worksheets.Item(2)

# Delete the entire 8th row (used to hold "t" in columns A-D).
$wsMember.Rows.Item(8).Delete()

# Row 7, column A changes from "r" to "e" (B7:D7 remain "r").
$wsMember.Range("A7").Value = "e"

# Replace the self-introduction text on the second sheet.
$newIntro = "asdfdsafasdfjlk`nasdjfjklsdajckldas`ndsanvkfdsaklcasdk`n"
$wsIntro.Range("A1").Value = $newIntro

# Narrow column A on the second sheet from ~21.17 to ~17.17 characters wide.
$wsIntro.Columns.Item(1).ColumnWidth = 16.333333333333332
